$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" (Total) sheet: shift the quarter rows (B:D, rows 2-9) down
#    by one row to make room for the new "2022-Q4" entry, then fill
#    in the new row 2 and the freshly-uncovered row 10. Column A (the
#    0-based running index) is left untouched for rows 2-9 and simply
#    extended with the next index (8) for the new row 10.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2:D9").Copy()
$total.Range("B3:D10").PasteSpecial(-4104)  # xlPasteAll

$total.Range("A2").Copy()
$total.Range("A10").PasteSpecial(-4122)     # xlPasteFormats
$total.Range("A10").Value = 8

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# ------------------------------------------------------------------
# 2) Insert a new "2022-Q4" sheet right after "总计" (i.e. before the
#    existing "2022-Q3" sheet) by duplicating "2022-Q3" (same layout:
#    one fund row) and overwriting its figures with the new quarter's
#    numbers.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $q3
$q4.Name = "2022-Q4"

# Columns D-G hold text-formatted numbers (e.g. "92.90"); force text
# number-format first so PowerShell/COM does not silently coerce the
# assigned string into a float (which would also eat trailing zeros).
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.43"
$q4.Range("E2").Value = "92.90"
$q4.Range("F2").Value = "2.89"
$q4.Range("G2").Value = "0.0124"
$q4.Range("H2").Value = 5

# ------------------------------------------------------------------
# 3) Keep "2020-Q4" as the active/selected tab, same as before the
#    edit (inserting/copying sheets shifts Excel's active tab, so
#    restore it explicitly).
# ------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
